$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Name row: set the Value column (B4) to "RoleVs"
$ws.Range("B4").Value = "RoleVs"

# Date row: update the Value column (B8) to the new generation timestamp
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
